# Weekly data refresh: a new weekly price record is inserted as row 388
# of the "Femacal de La Calera - Acelga" data sheet, pushing the existing
# rows 388:413 down to 389:414 (dimension grows from A1:R413 to A1:R414).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 388; everything at/after row 388 shifts down one.
$ws.Rows.Item(388).Insert()

# Populate the newly inserted row 388 with the new weekly record.
$ws.Cells.Item(388, 1).Value = 3
$ws.Cells.Item(388, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(388, 3).Value = "Coquimbo"
$ws.Cells.Item(388, 4).Value = 44826
$ws.Cells.Item(388, 5).Value = 5
$ws.Cells.Item(388, 6).Value = 100112009
$ws.Cells.Item(388, 7).Value = "Acelga"
$ws.Cells.Item(388, 8).Value = "Sin especificar"
$ws.Cells.Item(388, 9).Value = "Primera"
$ws.Cells.Item(388, 10).Value = 230
$ws.Cells.Item(388, 11).Value = 2500
$ws.Cells.Item(388, 12).Value = 2800
$ws.Cells.Item(388, 13).Value = 2657
$ws.Cells.Item(388, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(388, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(388, 16).Value = 443
$ws.Cells.Item(388, 17).Value = 6
$ws.Cells.Item(388, 18).Value = "Hortaliza"
